# Fantasy Football standings workbook edit
# - Clean out real names; use abbreviated team names for the draft board
#     "Aida's Astounding Team" -> "A's Astounding Team"
#     "Kelly's Deluxe Team"    -> "K's Deluxe Team"
#     "Magic Mikaela"          -> "Magic M"
# - These names appear many times in column A (one row per Team/Week combo)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCell = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162)
$lastRow = $lastCell.Row
if ($lastRow -lt 2) {
    $lastRow = 205
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq "Aida's Astounding Team") {
        $cell.Value2 = "A's Astounding Team"
    }
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq "Kelly's Deluxe Team") {
        $cell.Value2 = "K's Deluxe Team"
    }
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq "Magic Mikaela") {
        $cell.Value2 = "Magic M"
    }
}
